# New: better curated proteomics data
# Adds 13 new rows of curated proteomics-based kcat entries to the "F 0.5"
# worksheet (rows 20-32), matching rows already present (for the transport
# reactions and the two "priority 3" entries) on the "F 0.3" worksheet, plus
# two brand-new note strings describing proteomics-based kcat limits.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("F 0.5")

# Row 20
$ws.Cells.Item(20, 1).Value = "Q6C0K8"
$ws.Cells.Item(20, 2).Value = "YALI0F23837g"
$ws.Cells.Item(20, 3).Value = "YALI0F23837g"
$ws.Cells.Item(20, 4).Value = 2.643
$ws.Cells.Item(20, 4).NumberFormat = "0.00"
$ws.Cells.Item(20, 5).Value = "y000251"
$ws.Cells.Item(20, 6).Value = "Limits model when adding proteomics data. Calculated from specific activity of E. coli (EC 2.7.8.B10)"
$ws.Cells.Item(20, 7).Value = 1

# Row 21
$ws.Cells.Item(21, 1).Value = "F2Z694"
$ws.Cells.Item(21, 4).Value = 3054.42
$ws.Cells.Item(21, 6).Value = "Transport reaction. Kcat set to 10^3 order of magnitude"
$ws.Cells.Item(21, 7).Value = 1

# Row 22
$ws.Cells.Item(22, 1).Value = "Q6C1W2"
$ws.Cells.Item(22, 4).Value = 1930
$ws.Cells.Item(22, 6).Value = "Transport reaction. Kcat set to 10^3 order of magnitude"
$ws.Cells.Item(22, 7).Value = 1

# Row 23
$ws.Cells.Item(23, 1).Value = "Q6C0B0"
$ws.Cells.Item(23, 4).Value = 7992
$ws.Cells.Item(23, 6).Value = "Transport reaction. Kcat set to 10^3 order of magnitude"
$ws.Cells.Item(23, 7).Value = 1

# Row 24
$ws.Cells.Item(24, 1).Value = "Q6C428"
$ws.Cells.Item(24, 4).Value = 1216.9
$ws.Cells.Item(24, 6).Value = "Transport reaction. Kcat set to 10^3 order of magnitude"
$ws.Cells.Item(24, 7).Value = 1

# Row 25
$ws.Cells.Item(25, 1).Value = "Q6C7R0"
$ws.Cells.Item(25, 4).Value = 2081.9
$ws.Cells.Item(25, 6).Value = "Transport reaction. Kcat set to 10^3 order of magnitude"
$ws.Cells.Item(25, 7).Value = 1

# Row 26
$ws.Cells.Item(26, 1).Value = "Q6CCX5"
$ws.Cells.Item(26, 4).Value = 1930
$ws.Cells.Item(26, 6).Value = "Transport reaction. Kcat set to 10^3 order of magnitude"
$ws.Cells.Item(26, 7).Value = 1

# Row 27
$ws.Cells.Item(27, 1).Value = "Q6C3A8"
$ws.Cells.Item(27, 4).Value = 1930
$ws.Cells.Item(27, 6).Value = "Transport reaction. Kcat set to 10^3 order of magnitude"
$ws.Cells.Item(27, 7).Value = 1

# Row 28
$ws.Cells.Item(28, 1).Value = "Q6CAH9"
$ws.Cells.Item(28, 4).Value = 1249.6
$ws.Cells.Item(28, 6).Value = "Transport reaction. Kcat set to 10^3 order of magnitude"
$ws.Cells.Item(28, 7).Value = 1

# Row 29
$ws.Cells.Item(29, 1).Value = "Q6CG86"
$ws.Cells.Item(29, 4).Value = 7204.3
$ws.Cells.Item(29, 6).Value = "Transport reaction. Kcat set to 10^3 order of magnitude"
$ws.Cells.Item(29, 7).Value = 1

# Row 30
$ws.Cells.Item(30, 1).Value = "Q6C8F4"
$ws.Cells.Item(30, 4).Value = 1610.991
$ws.Cells.Item(30, 6).Value = "Transport reaction. Kcat set to 10^3 order of magnitude"
$ws.Cells.Item(30, 7).Value = 1

# Row 31
$ws.Cells.Item(31, 1).Value = "Q6CD72 + Q6C3F1"
$ws.Cells.Item(31, 4).Value = 4492
$ws.Cells.Item(31, 6).Value = "Transport reaction. Kcat set to 10^3 order of magnitude"
$ws.Cells.Item(31, 7).Value = "1 + 1"
$ws.Cells.Item(31, 7).HorizontalAlignment = -4152

# Row 32
$ws.Cells.Item(32, 1).Value = "Q6C8F2"
$ws.Cells.Item(32, 2).Value = "YALI0D20152g"
$ws.Cells.Item(32, 3).Value = "YALI0D20152g"
$ws.Cells.Item(32, 4).Value = 38
$ws.Cells.Item(32, 5).Value = "y000760"
$ws.Cells.Item(32, 6).Value = "Limits model when adding proteomics data. Using kcat of Aspergillus fumigatus (EC 2.3.1.4)"
$ws.Cells.Item(32, 7).Value = 1

# Update the active-cell selections to mirror the author's final view state.
$ws1 = $wb.Worksheets.Item("F 0.3")
$null = $ws1.Range("A40:G40").Select()
$null = $ws.Range("F30").Select()
